$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.932.81'
$ws.Range('E2').Value = '  +1.45%  '

$ws.Range('D3').Value = '1.670.88'
$ws.Range('E3').Value = '  +2.69%  '

$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.88'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +1.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.531'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  +5.92%  '

$ws.Range('E7').Value = '  +0.19%  '

$ws.Range('E8').Value = '  +2.61%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0620'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  +1.68%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.18'
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  +4.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').Style = $ws.Range('B11').Style
$ws.Range('E11').Value = '  +3.94%  '

$ws.Range('D12').Value = '1.908.53'
$ws.Range('E12').Value = '  +2.90%  '

$ws.Range('D13').Value = '1.663.92'
$ws.Range('E13').Value = '  +1.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.521'
$ws.Range('D15').Style = $ws.Range('B15').Style
$ws.Range('E15').Value = '  +1.67%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.57'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  +2.66%  '

$ws.Range('D17').Value = '26.959.71'
$ws.Range('E17').Value = '  +1.48%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.42'
$ws.Range('D18').Style = $ws.Range('B18').Style
$ws.Range('E18').Value = '  -0.84%  '

$ws.Range('D19').Value = '0.0₃0734'
$ws.Range('E19').Value = '  +1.22%  '

$ws.Range('E20').Value = '  -0.48%  '

$ws.Range('E21').Value = '  +0.21%  '

$ws.Range('E22').Value = '  +2.79%  '

$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.21'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +0.27%  '

$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.20'
$ws.Range('D24').Style = $ws.Range('B24').Style
$ws.Range('E24').Value = '  +0.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.59'
$ws.Range('D25').Style = $ws.Range('B25').Style
$ws.Range('E25').Value = '  -0.33%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.14'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  +0.92%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.116'
$ws.Range('D27').Style = $ws.Range('B27').Style
$ws.Range('E27').Value = '  +2.67%  '

$ws.Range('E28').Value = '  +1.06%  '

$ws.Range('E29').Value = '  +0.19%  '

$ws.Range('E30').Value = '  +0.54%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.17'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  +1.14%  '

$ws.Range('E32').Value = '  +1.59%  '

$ws.Range('D33').Value = '1.452.16'
$ws.Range('E33').Value = '  -4.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.15'
$ws.Range('D34').Style = $ws.Range('B34').Style
$ws.Range('E34').Value = '  +4.55%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.61'
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  +4.85%  '

$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.899'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +7.47%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.564'
$ws.Range('D38').Style = $ws.Range('B38').Style
$ws.Range('E38').Value = '  -0.83%  '

$ws.Range('E39').Value = '  +1.26%  '

$ws.Range('E40').Value = '  +3.18%  '

$ws.Range('E41').Value = '  +0.20%  '

$ws.Range('E42').Value = '  +4.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.74'
$ws.Range('D43').Style = $ws.Range('B43').Style
$ws.Range('E43').Value = '  +4.35%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.968'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  +6.61%  '

$ws.Range('D45').Value = '1.813.68'
$ws.Range('E45').Value = '  +2.64%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.779'
$ws.Range('D46').Style = $ws.Range('B46').Style
$ws.Range('E46').Value = '  +2.17%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.59'
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  +0.97%  '

$ws.Range('E48').Value = '  +0.82%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1000'
$ws.Range('D49').Style = $ws.Range('B49').Style
$ws.Range('E49').Value = '  +3.68%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0508'
$ws.Range('D50').Style = $ws.Range('B50').Style
$ws.Range('E50').Value = '  +1.30%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.60'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  +0.76%  '
